$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Plain text replacements (no risk of Excel auto-typing as a number/date)
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

$ws.Range("B15").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C15").Value = "519033 - Carlos Yujiro Shigue"

$ws.Range("B18").Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("C18").Value = "7290967 - Emerson Gonçalves de Melo"

# "01/01/2023" looks like a date to Excel's auto-typing, so setting it
# directly on B13/C13 would convert the cell to a date serial number and
# (if forced text via NumberFormat) would also mutate those cells' style
# index. Instead, stage the text in a scratch cell formatted as Text, copy
# it, and paste *values only* into the target cells so their original
# style is preserved untouched.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "01/01/2023"
$scratch.Copy()
$ws.Range("B13").PasteSpecial(-4163)
$ws.Range("C13").PasteSpecial(-4163)

# Remove the scratch column so dimensions/used-range match the original.
$ws.Columns.Item(26).Delete()
$excel.CutCopyMode = $false
